$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("cb155049-a59d-4c91-9454-a2fbac1c94a4", "Bawang", 210.0, "2023-01-16T09:00:57.893924Z", "2023-01-16T09:00:57.893924Z"),
    @("00003944-c211-4828-84e0-8db3ceaa60fb", "Brokoli", 310.0, "2023-01-16T09:01:05.133635Z", "2023-01-16T09:01:05.133635Z"),
    @("be6377de-e61e-4281-a6ea-9d20a5de0293", "Wortel", 410.0, "2023-01-16T09:01:11.259338Z", "2023-01-16T09:01:11.259338Z"),
    @("7804bc9d-17f7-4293-9575-4fbd60f0fb8c", "Jamur Trufle", 45.0, "2023-01-16T09:01:19.391102Z", "2023-01-16T09:01:19.391102Z")
)

$row = 3
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
